$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" updates ---
$ws = $wb.Worksheets.Item("Metadata")

# Publisher value (row 9, column B)
$ws.Range("B9").Value = "Independent Trusted Third Party of the University Medicine Greifswald"

# Contact value (row 10, column B)
$ws.Range("B10").Value = "Independent Trusted Third Party of the University Medicine Greifswald (https://www.ths-greifswald.de/)"

# Description value (row 12, column B)
$ws.Range("B12").Value = "Possible save actions in the context of adding patient identities. "
